# Auto-generated PowerShell script to update cryptos.xlsx per commit diff
# Commit message: Updated cryptos list on Sun Jul 28 14:24:53 UTC 2024 with GitHub Actions

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue 2 4 '67.806.73'
Set-TextValue 2 5 '  -1.80%  '

# Row 3
Set-TextValue 3 4 '3.265.54'
Set-TextValue 3 5 '  -1.25%  '

# Row 4
Set-TextValue 4 5 '  -0.01%  '

# Row 5
Set-TextValue 5 2 'BNB'
Set-TextValue 5 3 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
Set-TextValue 5 4 '580.28'
Set-TextValue 5 5 '  -1.50%  '

# Row 6
Set-TextValue 6 2 'Solana'
Set-TextValue 6 3 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue 6 4 '185.16'
Set-TextValue 6 5 '  -0.44%  '

# Row 7
Set-TextValue 7 5 '  +0.02%  '

# Row 8
Set-TextValue 8 4 '0.598'
Set-TextValue 8 5 '  -0.87%  '

# Row 9
Set-TextValue 9 2 'Dogecoin'
Set-TextValue 9 3 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue 9 4 '0.130'
Set-TextValue 9 5 '  -4.68%  '

# Row 10
Set-TextValue 10 2 'Toncoin'
Set-TextValue 10 3 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 10 4 '6.52'
Set-TextValue 10 5 '  -2.93%  '

# Row 11
Set-TextValue 11 2 'Cardano'
Set-TextValue 11 3 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextValue 11 4 '0.408'
Set-TextValue 11 5 '  -3.58%  '

# Row 12
Set-TextValue 12 2 'WrappedliquidstakedEther2.0'
Set-TextValue 12 3 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue 12 4 '3.829.72'
Set-TextValue 12 5 '  -1.44%  '

# Row 13
Set-TextValue 13 2 'TRON'
Set-TextValue 13 3 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue 13 4 '0.137'
Set-TextValue 13 5 '  +0.04%  '

# Row 14
Set-TextValue 14 2 'Avalanche'
Set-TextValue 14 3 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue 14 4 '27.43'
Set-TextValue 14 5 '  -6.49%  '

# Row 15
Set-TextValue 15 2 'WrappedBTC'
Set-TextValue 15 3 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue 15 4 '67.870.25'
Set-TextValue 15 5 '  -1.81%  '

# Row 16
Set-TextValue 16 2 'ShibaInu'
Set-TextValue 16 3 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue 16 4 '0.0000168'
Set-TextValue 16 5 '  -3.44%  '

# Row 17
Set-TextValue 17 2 'WrappedEther'
Set-TextValue 17 3 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 17 4 '3.269.80'
Set-TextValue 17 5 '  -0.44%  '

# Row 18
Set-TextValue 18 2 'Polkadot'
Set-TextValue 18 3 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 18 4 '5.70'
Set-TextValue 18 5 '  -3.21%  '

# Row 19
Set-TextValue 19 2 'Chainlink'
Set-TextValue 19 3 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 19 4 '13.50'
Set-TextValue 19 5 '  -1.95%  '

# Row 20
Set-TextValue 20 2 'BitcoinCash'
Set-TextValue 20 3 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 20 4 '399.41'
Set-TextValue 20 5 '  +2.49%  '

# Row 21
Set-TextValue 21 2 'Uniswap'
Set-TextValue 21 3 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue 21 4 '7.59'
Set-TextValue 21 5 '  -2.61%  '

# Row 22
Set-TextValue 22 2 'Dai'
Set-TextValue 22 3 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 22 4 '0.999'
Set-TextValue 22 5 '  +0.19%  '

# Row 23
Set-TextValue 23 2 'Litecoin'
Set-TextValue 23 3 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue 23 4 '71.10'
Set-TextValue 23 5 '  -1.39%  '

# Row 24
Set-TextValue 24 2 'Polygon'
Set-TextValue 24 3 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 24 4 '0.510'
Set-TextValue 24 5 '  -1.85%  '

# Row 25
Set-TextValue 25 2 'PEPE'
Set-TextValue 25 3 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue 25 4 '0.0000117'
Set-TextValue 25 5 '  -5.00%  '

# Row 26
Set-TextValue 26 2 'Kaspa'
Set-TextValue 26 3 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 26 4 '0.187'
Set-TextValue 26 5 '  -0.54%  '

# Row 27
Set-TextValue 27 2 'InternetComputer(DFINITY)'
Set-TextValue 27 3 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 27 4 '9.49'
Set-TextValue 27 5 '  -3.15%  '

# Row 28
Set-TextValue 28 2 'Binance-PegBSC-USD'
Set-TextValue 28 3 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue 28 4 '1.00'
Set-TextValue 28 5 '  +0.92%  '

# Row 29
Set-TextValue 29 2 'PancakeSwap'
Set-TextValue 29 3 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 29 4 '1.95'
Set-TextValue 29 5 '  -2.50%  '

# Row 30
Set-TextValue 30 2 'EthereumClassic'
Set-TextValue 30 3 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 30 4 '22.61'
Set-TextValue 30 5 '  -2.29%  '

# Row 31
Set-TextValue 31 2 'NEARProtocol'
Set-TextValue 31 3 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 31 4 '5.47'
Set-TextValue 31 5 '  -6.22%  '

# Row 32
Set-TextValue 32 2 'Aptos'
Set-TextValue 32 3 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 32 4 '6.92'
Set-TextValue 32 5 '  -3.72%  '

# Row 33
Set-TextValue 33 2 'USDe'
Set-TextValue 33 3 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue 33 4 '0.998'
Set-TextValue 33 5 '  +0.03%  '

# Row 34
Set-TextValue 34 4 '1.25'
Set-TextValue 34 5 '  -5.37%  '

# Row 35
Set-TextValue 35 2 'Monero'
Set-TextValue 35 3 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 35 4 '162.30'
Set-TextValue 35 5 '  -0.71%  '

# Row 36
Set-TextValue 36 2 'ImmutableX'
Set-TextValue 36 3 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 36 4 '1.46'
Set-TextValue 36 5 '  -6.29%  '

# Row 37
Set-TextValue 37 2 'Stacks'
Set-TextValue 37 3 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 37 4 '1.88'
Set-TextValue 37 5 '  -2.18%  '

# Row 38
Set-TextValue 38 2 'EnergySwap'
Set-TextValue 38 3 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 38 4 '26.78'
Set-TextValue 38 5 '  -0.08%  '

# Row 39
Set-TextValue 39 2 'Mantle'
Set-TextValue 39 3 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 39 4 '0.807'
Set-TextValue 39 5 '  -4.25%  '

# Row 40
Set-TextValue 40 2 'Filecoin'
Set-TextValue 40 3 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 40 4 '4.51'
Set-TextValue 40 5 '  -2.73%  '

# Row 41
Set-TextValue 41 2 'RenderToken'
Set-TextValue 41 3 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 41 4 '6.33'
Set-TextValue 41 5 '  -5.61%  '

# Row 42
Set-TextValue 42 2 'Maker'
Set-TextValue 42 3 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 42 4 '2.668.41'
Set-TextValue 42 5 '  +0.52%  '

# Row 43
Set-TextValue 43 2 'Hedera'
Set-TextValue 43 3 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 43 4 '0.0681'
Set-TextValue 43 5 '  -2.30%  '

# Row 44
Set-TextValue 44 2 'OKB'
Set-TextValue 44 3 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 44 4 '40.66'
Set-TextValue 44 5 '  -2.56%  '

# Row 45
Set-TextValue 45 2 'dogwifhat'
Set-TextValue 45 3 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 45 4 '2.42'
Set-TextValue 45 5 '  -8.48%  '

# Row 46
Set-TextValue 46 2 'Bittensor'
Set-TextValue 46 3 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 46 4 '334.77'
Set-TextValue 46 5 '  -1.58%  '

# Row 47
Set-TextValue 47 4 '24.55'
Set-TextValue 47 5 '  -3.89%  '

# Row 48
Set-TextValue 48 2 'VeChain'
Set-TextValue 48 3 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 48 4 '0.0275'
Set-TextValue 48 5 '  -4.08%  '

# Row 49
Set-TextValue 49 2 'Cosmos'
Set-TextValue 49 3 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 49 4 '6.36'
Set-TextValue 49 5 '  +0.39%  '

# Row 50
Set-TextValue 50 2 'Stellar'
Set-TextValue 50 3 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 50 4 '0.101'
Set-TextValue 50 5 '  -1.51%  '

# Row 51
Set-TextValue 51 2 'ONDO'
Set-TextValue 51 3 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-TextValue 51 4 '0.968'
Set-TextValue 51 5 '  -3.50%  '
